# Natmi following Dr Hou advice
# Recompute the C3/Cr2 ligand-receptor edge table across all three
# sending/target cluster combinations (ECs, FAPs, sCs), expanding the
# sheet from 4 data rows to the full 3x3 = 9 row matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C3"
$ws.Range("C2").Value = "Cr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 34.739995
$ws.Range("H2").Value = 104.219985
$ws.Range("I2").Value = 0.1827267341390226
$ws.Range("J2").Value = 0.1827267341390226
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.072165666666667
$ws.Range("N2").Value = 3.216497
$ws.Range("O2").Value = 0.6730206656141882
$ws.Range("P2").Value = 0.673020665614188
$ws.Range("Q2").Value = 37.24702989917166
$ws.Range("R2").Value = 335.223269092545
$ws.Range("S2").Value = 0.1229788682357518
$ws.Range("T2").Value = 0.1229788682357518

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C3"
$ws.Range("C3").Value = "Cr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 34.739995
$ws.Range("H3").Value = 104.219985
$ws.Range("I3").Value = 0.1827267341390226
$ws.Range("J3").Value = 0.1827267341390226
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1198123333333333
$ws.Range("N3").Value = 0.359437
$ws.Range("O3").Value = 0.0752086910034012
$ws.Range("P3").Value = 0.07520869100340119
$ws.Range("Q3").Value = 4.162279860938334
$ws.Range("R3").Value = 37.460518748445
$ws.Range("S3").Value = 0.01374263848592239
$ws.Range("T3").Value = 0.01374263848592239

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "C3"
$ws.Range("C4").Value = "Cr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 34.739995
$ws.Range("H4").Value = 104.219985
$ws.Range("I4").Value = 0.1827267341390226
$ws.Range("J4").Value = 0.1827267341390226
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.401087
$ws.Range("N4").Value = 1.203261
$ws.Range("O4").Value = 0.2517706433824107
$ws.Range("P4").Value = 0.2517706433824106
$ws.Range("Q4").Value = 13.933760374565
$ws.Range("R4").Value = 125.403843371085
$ws.Range("S4").Value = 0.04600522741734843
$ws.Range("T4").Value = 0.04600522741734842

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "C3"
$ws.Range("C5").Value = "Cr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 154.8642143333334
$ws.Range("H5").Value = 464.5926430000001
$ws.Range("I5").Value = 0.8145606273154508
$ws.Range("J5").Value = 0.8145606273154508
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.072165666666667
$ws.Range("N5").Value = 3.216497
$ws.Range("O5").Value = 0.6730206656141882
$ws.Range("P5").Value = 0.673020665614188
$ws.Range("Q5").Value = 166.0400936035079
$ws.Range("R5").Value = 1494.360842431571
$ws.Range("S5").Value = 0.5482161355789553
$ws.Range("T5").Value = 0.5482161355789552

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "C3"
$ws.Range("C6").Value = "Cr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 154.8642143333334
$ws.Range("H6").Value = 464.5926430000001
$ws.Range("I6").Value = 0.8145606273154508
$ws.Range("J6").Value = 0.8145606273154508
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1198123333333333
$ws.Range("N6").Value = 0.359437
$ws.Range("O6").Value = 0.0752086910034012
$ws.Range("P6").Value = 0.07520869100340119
$ws.Range("Q6").Value = 18.55464286911011
$ws.Range("R6").Value = 166.991785821991
$ws.Range("S6").Value = 0.06126203852330438
$ws.Range("T6").Value = 0.06126203852330437

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "C3"
$ws.Range("C7").Value = "Cr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 154.8642143333334
$ws.Range("H7").Value = 464.5926430000001
$ws.Range("I7").Value = 0.8145606273154508
$ws.Range("J7").Value = 0.8145606273154508
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.401087
$ws.Range("N7").Value = 1.203261
$ws.Range("O7").Value = 0.2517706433824107
$ws.Range("P7").Value = 0.2517706433824106
$ws.Range("Q7").Value = 62.11402313431368
$ws.Range("R7").Value = 559.0262082088232
$ws.Range("S7").Value = 0.2050824532131911
$ws.Range("T7").Value = 0.205082453213191

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "C3"
$ws.Range("C8").Value = "Cr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5157266666666667
$ws.Range("H8").Value = 1.54718
$ws.Range("I8").Value = 0.002712638545526686
$ws.Range("J8").Value = 0.002712638545526686
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.072165666666667
$ws.Range("N8").Value = 3.216497
$ws.Range("O8").Value = 0.6730206656141882
$ws.Range("P8").Value = 0.673020665614188
$ws.Range("Q8").Value = 0.5529444253844444
$ws.Range("R8").Value = 4.97649982846
$ws.Range("S8").Value = 0.001825661799481073
$ws.Range("T8").Value = 0.001825661799481073

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "C3"
$ws.Range("C9").Value = "Cr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5157266666666667
$ws.Range("H9").Value = 1.54718
$ws.Range("I9").Value = 0.002712638545526686
$ws.Range("J9").Value = 0.002712638545526686
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1198123333333333
$ws.Range("N9").Value = 0.359437
$ws.Range("O9").Value = 0.0752086910034012
$ws.Range("P9").Value = 0.07520869100340119
$ws.Range("Q9").Value = 0.06179041529555556
$ws.Range("R9").Value = 0.55611373766
$ws.Range("S9").Value = 0.0002040139941744322
$ws.Range("T9").Value = 0.0002040139941744322

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "C3"
$ws.Range("C10").Value = "Cr2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5157266666666667
$ws.Range("H10").Value = 1.54718
$ws.Range("I10").Value = 0.002712638545526686
$ws.Range("J10").Value = 0.002712638545526686
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.401087
$ws.Range("N10").Value = 1.203261
$ws.Range("O10").Value = 0.2517706433824107
$ws.Range("P10").Value = 0.2517706433824106
$ws.Range("Q10").Value = 0.2068512615533334
$ws.Range("R10").Value = 1.86166135398
$ws.Range("S10").Value = 0.0006829627518711804
$ws.Range("T10").Value = 0.0006829627518711803

